$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 3 "ECs" sending-cluster rows (old rows 2-4); remaining rows shift up
$ws.Rows("2:4").Delete()

# Update remaining rows (now rows 2-7, sending clusters FAPs/MuSCs) with refreshed TPM figures
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Col14a1"
$ws.Range("C2").Value = "Cd44"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 69.26966233333333
$ws.Range("H2").Value = 207.808987
$ws.Range("I2").Value = 0.9976788694238914
$ws.Range("J2").Value = 0.9976788694238914
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.142376000000001
$ws.Range("N2").Value = 24.427128
$ws.Range("O2").Value = 0.1741313933276368
$ws.Range("P2").Value = 0.1741313933276368
$ws.Range("Q2").Value = 564.0196361110374
$ws.Range("R2").Value = 5076.176724999336
$ws.Range("S2").Value = 0.1737272116263236
$ws.Range("T2").Value = 0.1737272116263236

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Col14a1"
$ws.Range("C3").Value = "Cd44"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 69.26966233333333
$ws.Range("H3").Value = 207.808987
$ws.Range("I3").Value = 0.9976788694238914
$ws.Range("J3").Value = 0.9976788694238914
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 24.34034433333333
$ws.Range("N3").Value = 73.021033
$ws.Range("O3").Value = 0.5205382400466131
$ws.Range("P3").Value = 0.5205382400466131
$ws.Range("Q3").Value = 1686.047433047063
$ws.Range("R3").Value = 15174.42689742357
$ws.Range("S3").Value = 0.5193300028216071
$ws.Range("T3").Value = 0.5193300028216071

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Col14a1"
$ws.Range("C4").Value = "Cd44"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 69.26966233333333
$ws.Range("H4").Value = 207.808987
$ws.Range("I4").Value = 0.9976788694238914
$ws.Range("J4").Value = 0.9976788694238914
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 14.277234
$ws.Range("N4").Value = 42.831702
$ws.Range("O4").Value = 0.3053303666257501
$ws.Range("P4").Value = 0.3053303666257501
$ws.Range("Q4").Value = 988.979178233986
$ws.Range("R4").Value = 8900.812604105875
$ws.Range("S4").Value = 0.3046216549759607
$ws.Range("T4").Value = 0.3046216549759607

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Col14a1"
$ws.Range("C5").Value = "Cd44"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.161158
$ws.Range("H5").Value = 0.483474
$ws.Range("I5").Value = 0.002321130576108561
$ws.Range("J5").Value = 0.002321130576108561
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.142376000000001
$ws.Range("N5").Value = 24.427128
$ws.Range("O5").Value = 0.1741313933276368
$ws.Range("P5").Value = 0.1741313933276368
$ws.Range("Q5").Value = 1.312209031408
$ws.Range("R5").Value = 11.809881282672
$ws.Range("S5").Value = 0.000404181701313164
$ws.Range("T5").Value = 0.0004041817013131639

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Col14a1"
$ws.Range("C6").Value = "Cd44"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.161158
$ws.Range("H6").Value = 0.483474
$ws.Range("I6").Value = 0.002321130576108561
$ws.Range("J6").Value = 0.002321130576108561
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 24.34034433333333
$ws.Range("N6").Value = 73.021033
$ws.Range("O6").Value = 0.5205382400466131
$ws.Range("P6").Value = 0.5205382400466131
$ws.Range("Q6").Value = 3.922641212071333
$ws.Range("R6").Value = 35.303770908642
$ws.Range("S6").Value = 0.001208237225005932
$ws.Range("T6").Value = 0.001208237225005931

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Col14a1"
$ws.Range("C7").Value = "Cd44"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.161158
$ws.Range("H7").Value = 0.483474
$ws.Range("I7").Value = 0.002321130576108561
$ws.Range("J7").Value = 0.002321130576108561
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 14.277234
$ws.Range("N7").Value = 42.831702
$ws.Range("O7").Value = 0.3053303666257501
$ws.Range("P7").Value = 0.3053303666257501
$ws.Range("Q7").Value = 2.300890476972
$ws.Range("R7").Value = 20.708014292748
$ws.Range("S7").Value = 0.0007087116497894656
$ws.Range("T7").Value = 0.0007087116497894655

